# Reorder the "Recorded By" (column G) contributor list on the
# "Session Analysis Results" sheet.
#
# Each G cell holds a comma-separated list of recorder names/emails
# (e.g. "System, dnasr281@gmail.com"). The last entry in the list is
# moved to the front, with the remaining entries keeping their
# original relative order (a right-rotation by one). Cells with a
# single name are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ",\s*"
    $count = $parts.Count

    if ($count -gt 1) {
        $newParts = @($parts[$count - 1]) + $parts[0..($count - 2)]
        $cell.Value = $newParts -join ", "
    }
}
